$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-8 to reflect regenerated
# strike-count-to-K values (regen save_data to use K instead of Strike#).
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 9
